$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price values look like plain numbers (e.g. "1.00", "0.140", "8.07").
# Excel would silently normalize these to numeric values/representations on
# assignment unless the cell is pre-formatted as Text, so force Text format on
# those specific cells first to preserve the exact source strings.
$priceTextCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D40", "D41", "D42", "D44", "D46", "D50", "D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "54.617.54"
$ws.Range("E2").Value = "  +5.51%  "
$ws.Range("D3").Value = "3.180.57"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "401.45"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("D6").Value = "109.89"
$ws.Range("E6").Value = "  +6.12%  "
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("D10").Value = "39.09"
$ws.Range("E10").Value = "  +4.43%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0892"
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.140"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "3.685.35"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "19.11"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "8.07"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("E16").Value = "  +7.45%  "
$ws.Range("D17").Value = "3.190.31"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "54.530.33"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("E21").Value = "  +4.39%  "
$ws.Range("D22").Value = "12.96"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("D23").Value = "72.28"
$ws.Range("E23").Value = "  +3.20%  "
$ws.Range("D24").Value = "276.33"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").Value = "3.25"
$ws.Range("E25").Value = "  +3.66%  "
$ws.Range("D26").Value = "8.05"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "27.82"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "7.58"
$ws.Range("E28").Value = "  +5.86%  "
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").Value = "11.03"
$ws.Range("E32").Value = "  +6.45%  "
$ws.Range("D33").Value = "0.0506"
$ws.Range("E33").Value = "  +12.65%  "
$ws.Range("D34").Value = "36.65"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "51.47"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  +6.49%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  +10.29%  "
$ws.Range("D40").Value = "4.09"
$ws.Range("E40").Value = "  +10.46%  "
$ws.Range("D41").Value = "1.93"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").Value = "0.292"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").Value = "131.64"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "22.17"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "2.095.18"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("B50").Value = "FlareNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr"
$ws.Range("D50").Value = "0.0513"
$ws.Range("E50").Value = "  +14.77%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0339"
$ws.Range("E51").Value = "  +6.35%  "
